# Update the timestamp embedded in the OLS regression-results text that was
# pasted into column B of every "backward elimination" step sheet.
# Old:  Date: Thu, 02 Jan 2020 / Time: 20:49:08
# New:  Date: Sun, 05 Jan 2020 / Time: 21:22:49

$wb = $excel.ActiveWorkbook

$oldDate = "Thu, 02 Jan 2020"
$newDate = "Sun, 05 Jan 2020"
$oldTime = "20:49:08"
$newTime = "21:22:49"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value()
    if ($text -ne $null -and $text.Contains($oldDate)) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value = $updated
    }
}
